$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column widths (best achievable values given COM's pixel-quantized width
# grid; chosen so the saved width matches the target as closely as possible)
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 23.666666666666664
$ws.Columns("B").ColumnWidth = 24.5
$ws.Columns("C").ColumnWidth = 30.333333333333332
$ws.Columns("D").ColumnWidth = 25.666666666666664
$ws.Columns("E").ColumnWidth = 38.16666666666667
$ws.Columns("F").ColumnWidth = 15.333333333333334
$ws.Columns("G").ColumnWidth = 24.333333333333332
$ws.Columns("H").ColumnWidth = 24.333333333333332
$ws.Columns("I").ColumnWidth = 16.666666666666664
$ws.Columns("J").ColumnWidth = 18.833333333333332

# ---------------------------------------------------------------------------
# E9 header cell gains wrap text (keeps its existing fill/border)
# ---------------------------------------------------------------------------
$ws.Range("E9").WrapText = $true

# ---------------------------------------------------------------------------
# Row 10 - first new test case (fully centered + wrap, taller row)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "TC_ADD CUSTOMER_001"
$ws.Range("B10").Value = "Verify the details of Customer while adding"
$ws.Range("C10").Value = "Enter Customer name,address ,valid phone no.,valid email ID ,valid Aadhar no.,but invalid PAN number"
$ws.Range("D10").Value = "Need a valid customer details to add customer"
$ws.Range("E10").Value = "Enter name ,address, mobile number,Aadhar number,PAN number and click on enter button."
$ws.Range("F10").Value = "<valid mobile no.> , <valid Aadharno>  <invalid PANno>"
$ws.Range("G10").Value = "A message ' PAN number is invalid or already exists' is shown"

$ws.Range("A10:G10").WrapText = $true
$ws.Range("A10:G10").VerticalAlignment = -4108
$ws.Range("C10").HorizontalAlignment = -4108

$ws.Rows("10").RowHeight = 85.5

# ---------------------------------------------------------------------------
# Row 11 - second new test case (wrap only, no centering)
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "TC_ADD CUSTOMER_001"
$ws.Range("B11").Value = "Verify the details of Customer while adding"
$ws.Range("C11").Value = "Enter Customer name,address ,valid phone no.,valid email ID ,invalid Aadhar no.,valid PAN number"
$ws.Range("D11").Value = "Need a valid customer details to add customer"
$ws.Range("E11").Value = "Enter name ,address, mobile number,Aadhar number,PAN number and click on enter button."
$ws.Range("F11").Value = "<valid mobile no.> <invalid Aadharno>  <valid PANno>"
$ws.Range("G11").Value = "A message ' Aadhar number is invalid or already exists' is shown"

$ws.Range("A11:G11").WrapText = $true

$ws.Rows("11").RowHeight = 61

# ---------------------------------------------------------------------------
# Row 12 - third new test case (wrap only, A12 left unstyled)
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "TC_ADD CUSTOMER_001"
$ws.Range("B12").Value = "Verify the details of Customer while adding"
$ws.Range("C12").Value = "Enter all other details valid but Date of birth invalid"
$ws.Range("D12").Value = "Need a valid customer details to add customer"
$ws.Range("E12").Value = "Enter name ,address, mobile number,Aadhar number,PAN number ,DOB and click on enter button."
$ws.Range("F12").Value = "<valid all_details>   <invalid DOB>"
$ws.Range("G12").Value = "A message 'DOB is invalid or customer below 18 years' is shown"

$ws.Range("B12:G12").WrapText = $true

$ws.Rows("12").RowHeight = 50

# ---------------------------------------------------------------------------
# Final selection / view state
# ---------------------------------------------------------------------------
$ws.Range("G12").Select() | Out-Null
